$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Test Case" header text (B1) ---
$ws.Range("B1").Value = "Test Case: Testing to see any errors while logging in"

# --- Rewrite step 1 (row 2) ---
$ws.Range("C2").Value = "Step 1: Sign-in with no fields filled with information"
$ws.Range("D2").Value = "An error text pops up saying that I need to fill out those field"

# --- Rewrite step 2 (row 3) ---
$ws.Range("C3").Value = "Step 2: Sign in with email filled out but no password"
$ws.Range("D3").Value = "An error text pops up saying that I need to fill out the password field"

# --- Rewrite step 3 (row 4) ---
$ws.Range("C4").Value = "Step 3: Sign in with password filled out but not email"
$ws.Range("D4").Value = "An error text pops up saying that the email is required"

# --- Add new step 4 (row 5) ---
$ws.Range("C5").Value = "Step 4: Fill out both fields with the wrong information"
$ws.Range("D5").Value = "An error pops up saying that the email or password is wrong"

# --- Add new step 5 (row 6) ---
$ws.Range("C6").Value = "Step 5: Fill out the fields with the correct information "
$ws.Range("D6").Value = "I am logged in and redirected to the user's dashboard with no errors."

# --- Update the selection / view state to match the new edited range ---
$ws.Range("B1:F6").Select() | Out-Null
